$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new label cells in row 8 ("L" / "C")
$ws.Range("C8").Value = "L"
$ws.Range("D8").Value = "C"

# Add the inductance value and the Lamor-frequency-derived capacitance formula in row 9
$ws.Range("C9").Value = 0.000000000417
$ws.Range("D9").Formula = '=1/(H6^2*C9*(2*PI())^2)'

# Both new numeric cells use scientific notation formatting
$ws.Range("C9:D9").NumberFormat = "0.00E+00"

# Move the active selection like in the final, saved worksheet state
[void]$ws.Range("D10").Select()
